$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '248.68'

Set-TextValue $ws.Range('D3') '21.72'

Set-TextValue $ws.Range('D4') '5.482'

Set-TextValue $ws.Range('D5') '0.05702'

Set-TextValue $ws.Range('D6') '3.369'

Set-TextValue $ws.Range('D7') '0.8035'

Set-TextValue $ws.Range('D8') '1.045'

Set-TextValue $ws.Range('D9') '0.1529'

Set-TextValue $ws.Range('D10') '0.07383'

Set-TextValue $ws.Range('D11') '0.03157'

Set-TextValue $ws.Range('D12') '0.03002'

Set-TextValue $ws.Range('D13') '0.09288'

Set-TextValue $ws.Range('B14') 'MCDex'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D14') '3.439'
Set-TextValue $ws.Range('E14') '13MCDexMCB'

Set-TextValue $ws.Range('B15') 'BitForexToken'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001641'
Set-TextValue $ws.Range('E15') '14BitForexTokenBF'

Set-TextValue $ws.Range('D16') '0.04714'

Set-TextValue $ws.Range('D17') '0.0005870'

Set-TextValue $ws.Range('D18') '0.006348'

Set-TextValue $ws.Range('D19') '0.005052'

Set-TextValue $ws.Range('D20') '0.001042'

Set-TextValue $ws.Range('D21') '0.0001499'

Set-TextValue $ws.Range('D22') '0.0003135'

Set-TextValue $ws.Range('D23') '3.769'

Set-TextValue $ws.Range('D24') '6.430'

Set-TextValue $ws.Range('D25') '2.113'

Set-TextValue $ws.Range('D26') '0.3284'

Set-TextValue $ws.Range('D40') '0.04111'

Set-TextValue $ws.Range('D41') '0.006957'

Set-TextValue $ws.Range('D42') '0.1046'

Set-TextValue $ws.Range('D43') '0.002969'

Set-TextValue $ws.Range('D44') '0.009127'

Set-TextValue $ws.Range('D45') '0.00005836'

Set-TextValue $ws.Range('D47') '0.0005500'

Set-TextValue $ws.Range('D48') '0.6824'

Set-TextValue $ws.Range('D49') '0.009298'

Set-TextValue $ws.Range('D50') '0.00002100'
